# Mise à jour de l'application
# Adds a new training-day column (BV) for 2025-11-04, fills in each
# player's attendance code for that day, and clears the trailing
# (no-longer-applicable) attendance cells for a player who left the team.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New date header in BV1 (04/11/2025, serial 45965) -----------------
# Set the value first (so the engine marks it / dependents dirty), then
# copy the number-format/style from BU1 so the cell keeps the same date
# style without touching the value we just set.
$ws.Range("BV1").Value = 45965
$ws.Range("BU1").Copy()
$ws.Range("BV1").PasteSpecial(-4122)   # xlPasteFormats

# --- 2. Attendance codes for the new day (column BV), rows 2-29 -----------
# Mapping of row -> attendance code for the 2025-11-04 session.
# Row 12, and the player in row 21 (who left before this date) are handled
# separately below.
$attendance = @{
    2  = "P";
    3  = "R";
    4  = "P";
    5  = "B";
    6  = "B";
    7  = "P";
    8  = "RH";
    9  = "M";
    10 = "P";
    11 = "P";
    13 = "B";
    14 = "P";
    15 = "P";
    16 = "B";
    17 = "P";
    18 = "B";
    19 = "P";
    20 = "P";
    22 = "P";
    23 = "RH";
    24 = "A";
    25 = "P";
    26 = "P";
    27 = "P";
    28 = "P";
    29 = "B"
}

foreach ($row in $attendance.Keys) {
    $target = "BV$row"
    $source = "BU$row"
    # Set the value first so the change is detected and dependent COUNTA /
    # COUNTIF formulas for that row recalculate, then copy the formatting
    # from the previous day's cell so the new cell matches its style.
    $ws.Range($target).Value = $attendance[$row]
    $ws.Range($source).Copy()
    $ws.Range($target).PasteSpecial(-4122)   # xlPasteFormats
}

# --- 3. Row 21: player left the team - clear trailing attendance cells ----
# BR21:BU21 no longer hold a value (kept their style only), and the new
# BV21 cell is created blank (style only, no value).
$ws.Range("BR21:BU21").ClearContents()
$ws.Range("BU21").Copy()
$ws.Range("BV21").PasteSpecial(-4122)   # xlPasteFormats

# --- 4. Update the selection shown in the sheet view ----------------------
$ws.Range("BX26").Select()

$excel.CutCopyMode = $false
